$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.061.98'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.652.26'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').Value = '218.42'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').Value = '0.5221'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').Value = '0.2618'
$ws.Range('E8').Value = '  -1.54%  '
$ws.Range('D9').Value = '0.06278'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = '20.50'
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('D11').Value = '0.07728'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '4.463'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '1.650.20'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.5436'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0₅8094'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '64.94'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.057.02'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = '1.003'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '4.571'
$ws.Range('E19').Value = '  -2.37%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '191.72'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '10.02'
$ws.Range('E21').Value = '  -1.46%  '
$ws.Range('B22').Value = 'Chainlink'
$ws.Range('C22').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D22').Value = '5.991'
$ws.Range('E22').Value = '  -2.41%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '138.59'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = '0.1232'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '7.256'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '16.13'
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.401'
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '0.05941'
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '1.279'
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '3.490'
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '3.233'
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = '1.542'
$ws.Range('E33').Value = '  -6.89%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '0.9467'
$ws.Range('E34').Value = '  -3.77%  '
$ws.Range('D35').Value = '2.412'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').Value = '2.755'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.5673'
$ws.Range('E37').Value = '  -4.79%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01602'
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '5.859'
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '0.8461'
$ws.Range('E40').Value = '  -2.10%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '100.63'
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.004.32'
$ws.Range('E43').Value = '  -4.37%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.794.86'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '56.65'
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '1.003'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.4298'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').Value = '7.918'
$ws.Range('E49').Value = '  -2.16%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05150'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').Value = '1.466'
$ws.Range('E51').Value = '  -0.72%  '
